# shading properties of Singapore archetypes
# - added shading type (type_shade) values to construction archetypes on the
#   ARCHITECTURE sheet, column L, rows 2:19 (previously blank cells become 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Fill in the new shading-type column (L) for every archetype row with 0
$ws.Range("L2:L19").Value = 0

# Make ARCHITECTURE the active sheet/tab and leave the new column selected,
# mirroring the author's final cursor position after entering the data
$ws.Activate()
$ws.Range("L2:L19").Select()
